# Refresh the crypto price/volume snapshot scraped by the scheduled job.
# Values are written as plain text to mirror the source cells (inlineStr
# in the sheet XML), so numeric-looking "Price" entries are pre-formatted
# as Text (via a single Union range -> one shared style) before being
# assigned; otherwise Excel would silently coerce strings such as
# "17.80" or "3.70" into the numbers 17.8 / 3.7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = $ws.Range("D5")
$textCells = $excel.Union($textCells, $ws.Range("D6"))
$textCells = $excel.Union($textCells, $ws.Range("D8"))
$textCells = $excel.Union($textCells, $ws.Range("D9"))
$textCells = $excel.Union($textCells, $ws.Range("D13"))
$textCells = $excel.Union($textCells, $ws.Range("D18"))
$textCells = $excel.Union($textCells, $ws.Range("D20"))
$textCells = $excel.Union($textCells, $ws.Range("D21"))
$textCells = $excel.Union($textCells, $ws.Range("D23"))
$textCells = $excel.Union($textCells, $ws.Range("D24"))
$textCells = $excel.Union($textCells, $ws.Range("D27"))
$textCells = $excel.Union($textCells, $ws.Range("D28"))
$textCells = $excel.Union($textCells, $ws.Range("D30"))
$textCells = $excel.Union($textCells, $ws.Range("D33"))
$textCells = $excel.Union($textCells, $ws.Range("D34"))
$textCells = $excel.Union($textCells, $ws.Range("D35"))
$textCells = $excel.Union($textCells, $ws.Range("D37"))
$textCells = $excel.Union($textCells, $ws.Range("D38"))
$textCells = $excel.Union($textCells, $ws.Range("D40"))
$textCells = $excel.Union($textCells, $ws.Range("D41"))
$textCells = $excel.Union($textCells, $ws.Range("D42"))
$textCells = $excel.Union($textCells, $ws.Range("D45"))
$textCells = $excel.Union($textCells, $ws.Range("D46"))
$textCells = $excel.Union($textCells, $ws.Range("D50"))
foreach ($area in $textCells.Areas) { $area.NumberFormat = "@" }

$updates = [ordered]@{
    "D2" = "54.176.36"
    "E2" = "  +0.42%  "
    "D3" = "2.264.66"
    "E3" = "  +1.21%  "
    "E4" = "  -0.51%  "
    "D5" = "496.38"
    "E5" = "  +0.35%  "
    "D6" = "128.77"
    "E6" = "  +1.10%  "
    "E7" = "  +0.03%  "
    "D8" = "0.525"
    "E8" = "  -0.42%  "
    "D9" = "0.0952"
    "E9" = "  +0.20%  "
    "E10" = "  +0.93%  "
    "E11" = "  +2.77%  "
    "E12" = "  +4.78%  "
    "D13" = "22.96"
    "E13" = "  +5.78%  "
    "D14" = "2.663.47"
    "E14" = "  +0.46%  "
    "D15" = "54.157.04"
    "E15" = "  +0.29%  "
    "E16" = "  +0.53%  "
    "D17" = "2.264.93"
    "E17" = "  +0.41%  "
    "D18" = "10.22"
    "E18" = "  +2.21%  "
    "E19" = "  +1.14%  "
    "D20" = "303.12"
    "E20" = "  +1.00%  "
    "D21" = "6.33"
    "E21" = "  -1.30%  "
    "E22" = "  +0.38%  "
    "D23" = "60.57"
    "E23" = "  -2.44%  "
    "D24" = "0.991"
    "E24" = "  -2.80%  "
    "E25" = "  +0.05%  "
    "E26" = "  +3.16%  "
    "D27" = "172.66"
    "E27" = "  +2.84%  "
    "D28" = "1.61"
    "E28" = "  -0.06%  "
    "B29" = "PEPE"
    "C29" = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
    "D29" = "0.0₃0690"
    "E29" = "  +0.86%  "
    "B30" = "Aptos"
    "C30" = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
    "D30" = "5.96"
    "E30" = "  +1.68%  "
    "E31" = "  +0.91%  "
    "E32" = "  -0.02%  "
    "D33" = "17.80"
    "E33" = "  +0.56%  "
    "D34" = "0.996"
    "E34" = "  +0.43%  "
    "D35" = "0.944"
    "E35" = "  +4.35%  "
    "E36" = "  +1.38%  "
    "D37" = "3.70"
    "E37" = "  +0.51%  "
    "D38" = "0.375"
    "E38" = "  +0.48%  "
    "E39" = "  +0.07%  "
    "D40" = "3.37"
    "E40" = "  +0.63%  "
    "D41" = "4.80"
    "E41" = "  -1.43%  "
    "D42" = "124.54"
    "E42" = "  -1.07%  "
    "E43" = "  +1.77%  "
    "E44" = "  +0.98%  "
    "D45" = "0.544"
    "E45" = "  +0.28%  "
    "D46" = "241.67"
    "E46" = "  +2.11%  "
    "E47" = "  +0.85%  "
    "E48" = "  +1.12%  "
    "E49" = "  +0.86%  "
    "D50" = "16.13"
    "E50" = "  -0.34%  "
    "E51" = "  -0.38%  "
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
